$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting the existing rows 94-96 down to 95-97.
$ws.Rows.Item(94).Insert()

# Fill the new row 94 with the new weekly data point (matches target diff).
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 44509
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 100112012
$ws.Range("G94").Value = "Espinaca"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 30
$ws.Range("K94").Value = 8000
$ws.Range("L94").Value = 8000
$ws.Range("M94").Value = 8000
$ws.Range("N94").Value = "$/docena de atados"
$ws.Range("O94").Value = "Región de La Araucanía"
$ws.Range("P94").Value = 2667
$ws.Range("Q94").Value = 3
$ws.Range("R94").Value = "Hortaliza"
